$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NewUser")

# Update the test-account row with the newly generated credentials.
$ws.Range("A2").Value = "TestPF1221+12082021182038@gmail.com"
$ws.Range("C2").Value = "TestPF1221_12082021182038"

# The new password is all digits ("12082021182154"); a plain .Value
# assignment would be auto-typed as a number, but the source cell holds
# text. Build it as a text-formula result, then copy/paste-values so it
# lands as a literal text value (matching the original cell's string
# type) without leaving behind any extra cell-style definition.
$ws.Range("B2").Formula = "=""12082021182154"""
$ws.Range("B2").Copy()
$ws.Range("B2").PasteSpecial(-4163)

# NewUser becomes the active sheet/selection after the edit.
$ws.Activate()
$ws.Range("B3").Select()
